# The two oldest measurements (rows 2-3) were dropped from the report and
# two new measurements were appended at the end, with the remaining rows
# shifting up to fill the gap (ids stay sequential: 10..19 instead of 8..17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the first two data rows; everything below shifts up by two.
$ws.Rows("2:3").Delete()

# Append the two new measurements captured for 2024-05-07 at the bottom.
$ws.Range("A10").Value = 45419.86436986816
$ws.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B10").Value = 3.1
$ws.Range("C10").Value = 1.39
$ws.Range("D10").Value = "WindDirection.NORTH_WEST"
$ws.Range("E10").Value = 989.7
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = "Сколково"
$ws.Range("J10").Value = 18

$ws.Range("A11").Value = 45419.86855657802
$ws.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 1.27
$ws.Range("D11").Value = "WindDirection.NORTH_WEST"
$ws.Range("E11").Value = 989.9
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = "Сколково"
$ws.Range("J11").Value = 19
